$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing rows 5-7: convert a few stray "0" text cells to genuine numeric 0 ---
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("H7").Value = 0

# --- Part 1: protect number/date-like text cells as Text before assigning ---
$textCells = @(
    "A19", "I19", "A20", "I20", "A21", "I21", "A22", "D22", "E22", "H22", "I22", "A23", "H23", "I23", "A24", "D24", "E24", "F24", "G24", "H24", "I24", "A25", "I25", "A26", "I26", "A27", "I27", "A28", "I28", "A29", "I29", "A30", "I30", "A31", "I31", "A32", "I32", "A33", "I33", "A34", "I34", "A35", "I35"
)
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# --- Part 2: assign cell values for new rows 19-35 ---
$ws.Range("A19").Value = "2025-10-20"
$ws.Range("B19").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice20102025.pdf"
$ws.Range("C19").Value = "CHANDERIYA `nLEAD ZINC `nSMELTER"
$ws.Range("E19").Value = "315,100 316,600 315,600 314,600 313,100"
$ws.Range("I19").Value = "203,800"
$ws.Range("A20").Value = "2025-10-20"
$ws.Range("B20").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice20102025.pdf"
$ws.Range("C20").Value = "HYDRO-1 UNIT"
$ws.Range("E20").Value = "315,100 316,600 315,600 314,600 313,100"
$ws.Range("I20").Value = "203,800"
$ws.Range("A21").Value = "2025-10-20"
$ws.Range("B21").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice20102025.pdf"
$ws.Range("C21").Value = "NEW HYDRO `nSMELTER `nCHANDERIYA"
$ws.Range("E21").Value = "315,100 316,600 315,600 314,600 313,100"
$ws.Range("I21").Value = "203,800"
$ws.Range("A22").Value = "2025-10-20"
$ws.Range("B22").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice20102025.pdf"
$ws.Range("C22").Value = "ZINC SMELTER `nDEBRI"
$ws.Range("D22").Value = "0"
$ws.Range("E22").Value = "0"
$ws.Range("G22").Value = "0  314,600"
$ws.Range("H22").Value = "0"
$ws.Range("I22").Value = "0"
$ws.Range("A23").Value = "2025-10-20"
$ws.Range("B23").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice20102025.pdf"
$ws.Range("C23").Value = "Pantnagar `nMelting&Castin `ngPlant"
$ws.Range("E23").Value = "315,100 316,600 315,600 314,600"
$ws.Range("H23").Value = "0"
$ws.Range("I23").Value = "203,800"
$ws.Range("A24").Value = "2025-10-20"
$ws.Range("B24").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice20102025.pdf"
$ws.Range("C24").Value = "RAJPURA DARIBA `nLEAD SMELTER"
$ws.Range("D24").Value = "0"
$ws.Range("E24").Value = "0"
$ws.Range("F24").Value = "0"
$ws.Range("G24").Value = "0"
$ws.Range("H24").Value = "0"
$ws.Range("I24").Value = "203,800"
$ws.Range("A25").Value = "2025-10-20"
$ws.Range("B25").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice20102025.pdf"
$ws.Range("C25").Value = "Faridabad `nDepot"
$ws.Range("E25").Value = "317,600 319,100 313,100 317,100 315,600"
$ws.Range("I25").Value = "206,300"
$ws.Range("A26").Value = "2025-10-20"
$ws.Range("B26").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice20102025.pdf"
$ws.Range("C26").Value = "Panvel Depot"
$ws.Range("E26").Value = "318,400 319,900 318,900 317,900 316,400"
$ws.Range("I26").Value = "206,700"
$ws.Range("A27").Value = "2025-10-20"
$ws.Range("B27").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice20102025.pdf"
$ws.Range("C27").Value = "Pune Depot"
$ws.Range("E27").Value = "318,400 319,900 318,900 317,900 316,400"
$ws.Range("I27").Value = "207,100"
$ws.Range("A28").Value = "2025-10-20"
$ws.Range("B28").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice20102025.pdf"
$ws.Range("C28").Value = "Baroda Depot"
$ws.Range("E28").Value = "318,400 319,900 318,900 317,900 316,400"
$ws.Range("I28").Value = "207,100"
$ws.Range("A29").Value = "2025-10-20"
$ws.Range("B29").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice20102025.pdf"
$ws.Range("C29").Value = "Raipur Depot"
$ws.Range("E29").Value = "318,400 319,900 318,900 317,900 316,400"
$ws.Range("I29").Value = "207,100"
$ws.Range("A30").Value = "2025-10-20"
$ws.Range("B30").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice20102025.pdf"
$ws.Range("C30").Value = "JAMSHEDPUR `nSTOCK POINT"
$ws.Range("E30").Value = "316,100 317,600 316,600 315,600 314,100"
$ws.Range("I30").Value = "204,800"
$ws.Range("A31").Value = "2025-10-20"
$ws.Range("B31").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice20102025.pdf"
$ws.Range("D31").Value = "Kolkata Depot  316,100 317,600 316,600 315,600 314,100"
$ws.Range("I31").Value = "204,800"
$ws.Range("A32").Value = "2025-10-20"
$ws.Range("B32").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice20102025.pdf"
$ws.Range("C32").Value = "Bangalore `nDepot"
$ws.Range("E32").Value = "316,100 317,600 316,600 315,600 314,100"
$ws.Range("I32").Value = "204,800"
$ws.Range("A33").Value = "2025-10-20"
$ws.Range("B33").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice20102025.pdf"
$ws.Range("C33").Value = "Hyderabad `nDepot"
$ws.Range("E33").Value = "316,100 317,600 316,600 315,600 314,100"
$ws.Range("I33").Value = "204,800"
$ws.Range("A34").Value = "2025-10-20"
$ws.Range("B34").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice20102025.pdf"
$ws.Range("D34").Value = "Chennai Depot  316,100 317,600 316,600 315,600 314,100"
$ws.Range("I34").Value = "204,800"
$ws.Range("A35").Value = "2025-10-20"
$ws.Range("B35").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice20102025.pdf"
$ws.Range("C35").Value = "Sindesar `nsmelter HZAPL"
$ws.Range("E35").Value = "315,100 316,600"
$ws.Range("G35").Value = "0  314,600 313,100"
$ws.Range("I35").Value = "203,800"

# --- Part 3: keep row heights at sheet default (cells with embedded newlines otherwise auto-grow) ---
$multilineRows = @(19, 21, 22, 23, 24, 25, 30, 32, 33, 35)
foreach ($r in $multilineRows) {
    $ws.Rows.Item($r).RowHeight = 15
}
